# Calibrate to Globocan 2018 CC inc by age
# - Stop calibrating to Globocan 2012 cervical cancer incidence data by age
#   (flip "Usage Status" flag from Y to N for the Globocan 2012 by-age rows)
# - Switch to calibrating to Globocan 2018 cervical cancer incidence data by age
#   (append a new block of rows, one per age group, pulled from Globocan 2018)
# - Load and save Globocan 2018 data (new shared strings / source labels)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Stop using the Globocan 2012 by-age incidence rows (275-287) for
#    calibration: "Usage Status" column L flips from "Y" to "N".
# ---------------------------------------------------------------------------
for ($r = 275; $r -le 287; $r++) {
    $ws.Cells.Item($r, 12).Value = "N"
}

# ---------------------------------------------------------------------------
# 2) Append the Globocan 2018 by-age incidence rows (294-305), mirroring the
#    structure of the Globocan 2012 block (274-287).
# ---------------------------------------------------------------------------
$ages = @("15–19 years","20–24 years","25–29 years","30–34 years","35–39 years","40–44 years","45–49 years","50–54 years","55–59 years","60–64 years","65-69 years","70-74 years")
$gvals = @(3.9070527772503256,19.718592555074853,35.815396540774046,53.753954608683522,71.181262777416052,85.090519192952186,95.443618498222065,96.099524161750054,94.851577973825016,97.48752644071628,101.98823806664386,110.42515650113373)

for ($i = 0; $i -lt 12; $i++) {
    $r = 294 + $i

    $ws.Cells.Item($r, 1).Value = "CC Incidence in Women"
    $ws.Cells.Item($r, 2).Value = "Globocan 2018"
    $ws.Cells.Item($r, 3).Value = $ages[$i]
    $ws.Cells.Item($r, 4).Value = 2018

    $g = $ws.Cells.Item($r, 7)
    $g.Value = $gvals[$i]
    $g.NumberFormat = "0.00"
    $g.HorizontalAlignment = -4131

    $h = $ws.Cells.Item($r, 8)
    $h.Formula = "=G" + $r
    $h.NumberFormat = "0.00"

    $ii = $ws.Cells.Item($r, 9)
    $ii.Formula = "=H" + $r
    $ii.NumberFormat = "0.00"

    $ws.Cells.Item($r, 12).Value = "Y"
}

# Row 294 gets a taller row (wraps the comment text) and the long Globocan
# 2018 methodology note in column M, formatted like the matching Globocan
# 2012 note already on row 274.
$ws.Rows.Item(294).RowHeight = 60

$note = $ws.Cells.Item(294, 13)
$note.Value = "Incidence rates per 100,000 women come from Globocan 2018 from Minttu. Assume a normal approximation of the Poisson distribution where µ=λ and variance=λ  and λ=CC incidence rate per year."
$note.WrapText = $true
$noteChars = $note.Characters(136, 53)
$noteChars.Font.Name = "Calibri"
$noteChars.Font.Size = 11

# ---------------------------------------------------------------------------
# 3) Spacer row 306 (blank, keeps the same column styling as the data block)
#    and the closing formatting-only row 310.
# ---------------------------------------------------------------------------
$ws.Cells.Item(306, 7).NumberFormat = "0.00"

$ws.Rows.Item(310).RowHeight = 15.75
for ($c = 3; $c -le 14; $c++) {
    $cell310 = $ws.Cells.Item(310, $c)
    $cell310.Font.Size = 12
}

# ---------------------------------------------------------------------------
# 4) View settings: zoom out a bit and leave the selection on the new last
#    row of data.
# ---------------------------------------------------------------------------
$aw = $excel.ActiveWindow
$aw.Zoom = 70
$ws.Range("E310").Select()
